$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (D) and Volume(1h) (E) updates ---
# Leading "'" forces Excel to treat the numeric-looking price strings as
# plain text (matching the source inlineStr cells); Style is reset to "Normal"
# afterwards so no stray text-format style is left on the cell.

$ws.Range("D2").Value = "'29.150.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "'1.841.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'244.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").Value = "'0.6257"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.24%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'0.07497"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("D9").Value = "'0.2938"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("D10").Value = "'23.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("D11").Value = "'0.07716"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").Value = "'1.878.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.34%  "
$ws.Range("D13").Value = "'5.021"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").Value = "'0.6759"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").Value = "'83.06"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").Value = "'0.000009288"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.27%  "
$ws.Range("D17").Value = "'5.970"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.05%  "
$ws.Range("D18").Value = "'29.157.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Value = "'2.130.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("D20").Value = "'230.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.80%  "
$ws.Range("D21").Value = "'12.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").Value = "'7.191"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").Value = "'1.002"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "'160.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").Value = "'17.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D29").Value = "'1.505"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").Value = "'4.191"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.56%  "
$ws.Range("D31").Value = "'4.158"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.08%  "
$ws.Range("D32").Value = "'0.05586"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.77%  "
$ws.Range("D33").Value = "'1.206"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").Value = "'0.7509"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'1.855"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.55%  "
$ws.Range("D36").Value = "'1.145"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").Value = "'2.661"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").Value = "'2.773"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("D39").Value = "'1.226.57"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.31%  "
$ws.Range("D40").Value = "'0.01788"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").Value = "'6.575"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("D42").Value = "'0.9027"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.46%  "
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").Value = "'2.020.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.74%  "
$ws.Range("D45").Value = "'102.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("D46").Value = "'66.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.23%  "
$ws.Range("D47").Value = "'0.00000000123"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("D48").Value = "'0.5104"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.4090"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").Value = "'9.165"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.02%  "
$ws.Range("D51").Value = "'0.05842"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.15%  "

# --- Row 26/27 swap: Cosmos <-> Stellar ---
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").Value = "'0.1394"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.56%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'8.556"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.06%  "
